$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 62, pushing the existing rows 62-63
# (Banquete / Primera, week of 2021-10-22, Región Metropolitana) down to
# rows 64-65 unchanged.
$ws.Range("A62:R63").EntireRow.Insert()

# New row 62: Banquete, week of 2021-11-09, origin Provincia de Linares
$ws.Cells.Item(62, 1).Value = 12
$ws.Cells.Item(62, 2).Value = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(62, 3).Value = "Metropolitana"
$ws.Cells.Item(62, 4).Value = 44509
$ws.Cells.Item(62, 5).Value = 13
$ws.Cells.Item(62, 6).Value = 300000000
$ws.Cells.Item(62, 7).Value = "Espárragos"
$ws.Cells.Item(62, 8).Value = "Sin especificar"
$ws.Cells.Item(62, 9).Value = "Banquete"
$ws.Cells.Item(62, 10).Value = 280
$ws.Cells.Item(62, 11).Value = 1200
$ws.Cells.Item(62, 12).Value = 1200
$ws.Cells.Item(62, 13).Value = 1200
$ws.Cells.Item(62, 14).Value = "$/kilo"
$ws.Cells.Item(62, 15).Value = "Provincia de Linares"
$ws.Cells.Item(62, 16).Value = 1200
$ws.Cells.Item(62, 17).Value = 1
$ws.Cells.Item(62, 18).Value = "Hortaliza"

# New row 63: Primera, week of 2021-11-09, origin Provincia de Linares
$ws.Cells.Item(63, 1).Value = 12
$ws.Cells.Item(63, 2).Value = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(63, 3).Value = "Metropolitana"
$ws.Cells.Item(63, 4).Value = 44509
$ws.Cells.Item(63, 5).Value = 13
$ws.Cells.Item(63, 6).Value = 300000000
$ws.Cells.Item(63, 7).Value = "Espárragos"
$ws.Cells.Item(63, 8).Value = "Sin especificar"
$ws.Cells.Item(63, 9).Value = "Primera"
$ws.Cells.Item(63, 10).Value = 250
$ws.Cells.Item(63, 11).Value = 1000
$ws.Cells.Item(63, 12).Value = 1000
$ws.Cells.Item(63, 13).Value = 1000
$ws.Cells.Item(63, 14).Value = "$/kilo"
$ws.Cells.Item(63, 15).Value = "Provincia de Linares"
$ws.Cells.Item(63, 16).Value = 1000
$ws.Cells.Item(63, 17).Value = 1
$ws.Cells.Item(63, 18).Value = "Hortaliza"

# Make sure the used range / dimension reflects the two new rows.
$ws.Range("D62:D63").NumberFormat = $ws.Range("D64").NumberFormat
